$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Fix the title: "Encuesta de evaluación de calida" + bookmark +
#    "d" -> single run "Encuesta de evaluación de calidad" (no
#    bookmark). Find/Replace across the whole text stream merges the
#    two runs and drops the now-redundant _GoBack bookmark that used
#    to split them.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Encuesta de evaluación de calidad", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Encuesta de evaluación de calidad", 2) | Out-Null

# ------------------------------------------------------------------
# Locate the quality-survey table (the only table in the document).
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)

# ------------------------------------------------------------------
# 2. Row "Herramientas proporcionadas para la gestión de reuniones"
#    becomes three runs: "...para la " / "planificación " / "de
#    reuniones", and loses its <w:pPr><w:jc w:val="both"/></w:pPr>.
#    Row index discovered via probing: table row 6 (1-based).
# ------------------------------------------------------------------
$xmlPlanificacion = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">Herramientas proporcionadas para la </w:t></w:r>
<w:r><w:t xml:space="preserve">planificación </w:t></w:r>
<w:r><w:t>de reuniones</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rowPlanificacion = $null
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    if ($t.Cell($r, 1).Range.Text -like "Herramientas proporcionadas para la gestión de reuniones*") {
        $rowPlanificacion = $r
        break
    }
}
$t.Cell($rowPlanificacion, 1).Range.InsertXML($xmlPlanificacion) | Out-Null

# ------------------------------------------------------------------
# 3. The empty row right below gets the new "comunicación en las
#    reuniones" line, again as three runs, again losing its pPr.
# ------------------------------------------------------------------
$xmlComunicacion = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">Herramientas proporcionadas para la </w:t></w:r>
<w:r><w:t xml:space="preserve">comunicación en las </w:t></w:r>
<w:r><w:t>reuniones</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$rowComunicacion = $rowPlanificacion + 1
$t.Cell($rowComunicacion, 1).Range.InsertXML($xmlComunicacion) | Out-Null

# ------------------------------------------------------------------
# 4. The next empty row keeps its pPr (jc=both) untouched, but now
#    carries the relocated _GoBack bookmark (collapsed, empty range).
# ------------------------------------------------------------------
$rowBookmark = $rowComunicacion + 1
$bmRange = $t.Cell($rowBookmark, 1).Range
$bmRange.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
